# Atualizado por script em 12-11-2023 14:45
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helpers (positional params only - named params silently fail to bind here) ---

function Get-Record {
    param($Row)
    $f = $ws.Range("F$Row").Value2
    $g = $ws.Range("G$Row").Value2
    $h = $ws.Range("H$Row").Value2
    $i = $ws.Range("I$Row").Value2
    $j = $ws.Range("J$Row").Value2
    $k = $ws.Range("K$Row").Value2
    $l = $ws.Range("L$Row").Value2
    $m = $ws.Range("M$Row").Value2
    $n = $ws.Range("N$Row").Value2
    $o = $ws.Range("O$Row").Value2
    $p = $ws.Range("P$Row").Value2
    $q = $ws.Range("Q$Row").Value2
    $r = $ws.Range("R$Row").Value2
    $s = $ws.Range("S$Row").Value2
    $t = $ws.Range("T$Row").Value2
    $u = $ws.Range("U$Row").Value2
    $v = $ws.Range("V$Row").Value2
    [PSCustomObject]@{
        F = $f; G = $g; H = $h; I = $i; J = $j; K = $k; L = $l; M = $m
        N = $n; O = $o; P = $p; Q = $q; R = $r; S = $s; T = $t; U = $u; V = $v
    }
}

function Set-Record {
    param($Row, $Rec)
    $ws.Range("F$Row").Value = $Rec.F
    $ws.Range("G$Row").Value = $Rec.G
    $ws.Range("H$Row").Value = $Rec.H
    $ws.Range("I$Row").Value = $Rec.I
    $ws.Range("J$Row").Value = $Rec.J
    $ws.Range("K$Row").Value = $Rec.K
    $ws.Range("L$Row").Value = $Rec.L
    $ws.Range("M$Row").Value = $Rec.M
    $ws.Range("N$Row").Value = $Rec.N
    $ws.Range("O$Row").Value = $Rec.O
    $ws.Range("P$Row").Value = $Rec.P
    $ws.Range("Q$Row").Value = $Rec.Q
    $ws.Range("R$Row").Value = $Rec.R
    $ws.Range("S$Row").Value = $Rec.S
    $ws.Range("T$Row").Value = $Rec.T
    $ws.Range("U$Row").Value = $Rec.U
    $ws.Range("V$Row").Value = $Rec.V
}

# ---------------------------------------------------------------------------
# Rows 8 and 9 swap their match data (F:V) - A:E (index/country/tourney/date)
# stay put on each physical row.
# ---------------------------------------------------------------------------
$r8 = Get-Record 8
$r9 = Get-Record 9
Set-Record 8 $r9
Set-Record 9 $r8

# ---------------------------------------------------------------------------
# Rows 74, 75, 76 rotate: new74 = old76, new75 = old74, new76 = old75
# ---------------------------------------------------------------------------
$r74 = Get-Record 74
$r75 = Get-Record 75
$r76 = Get-Record 76
Set-Record 74 $r76
Set-Record 75 $r74
Set-Record 76 $r75

# ---------------------------------------------------------------------------
# Rows 85 and 86 swap their match data (F:V)
# ---------------------------------------------------------------------------
$r85 = Get-Record 85
$r86 = Get-Record 86
Set-Record 85 $r86
Set-Record 86 $r85

# ---------------------------------------------------------------------------
# Rows 91 and 92 swap their match data (F:V)
# ---------------------------------------------------------------------------
$r91 = Get-Record 91
$r92 = Get-Record 92
Set-Record 91 $r92
Set-Record 92 $r91

# ---------------------------------------------------------------------------
# Append two new match rows (94 and 95), mirroring row 93's formatting for
# columns A (bold/centered index) and E (date-time number format).
# ---------------------------------------------------------------------------
$ws.Range("A93").Copy()
$ws.Range("A94").PasteSpecial(-4122)
$ws.Range("A93").Copy()
$ws.Range("A95").PasteSpecial(-4122)
$ws.Range("E93").Copy()
$ws.Range("E94").PasteSpecial(-4122)
$ws.Range("E93").Copy()
$ws.Range("E95").PasteSpecial(-4122)

$ws.Range("A94").Value = 93
$ws.Range("B94").Value = "denmark"
$ws.Range("C94").Value = "1st-division"
$ws.Range("D94").Value = "2023-2024"
$ws.Range("E94").Value = 45242.54166666666
$rec94 = [PSCustomObject]@{
    F = "Naestved"; G = 1; H = "Aalborg"; I = 1; J = 5.5; K = "06/11/2023 21:12"
    L = 5.23; M = "12/11/2023 12:58"; N = 4.26; O = "06/11/2023 21:12"
    P = 4.15; Q = "12/11/2023 12:58"; R = 1.56; S = "06/11/2023 21:12"
    T = 1.62; U = "12/11/2023 12:58"
    V = "https://www.betexplorer.com/football/denmark/1st-division/naestved-if-aalborg/0Gbq0OCb/"
}
Set-Record 94 $rec94

$ws.Range("A95").Value = 94
$ws.Range("B95").Value = "denmark"
$ws.Range("C95").Value = "1st-division"
$ws.Range("D95").Value = "2023-2024"
$ws.Range("E95").Value = 45242.58333333334
$rec95 = [PSCustomObject]@{
    F = "Kolding IF"; G = 0; H = "Vendsyssel"; I = 1; J = 1.92; K = "05/11/2023 14:12"
    L = 1.81; M = "12/11/2023 13:59"; N = 3.76; O = "05/11/2023 14:12"
    P = 3.51; Q = "12/11/2023 13:59"; R = 3.72; S = "05/11/2023 14:12"
    T = 4.79; U = "12/11/2023 13:59"
    V = "https://www.betexplorer.com/football/denmark/1st-division/kolding-if-vendsyssel-ff/lxVQsscU/"
}
Set-Record 95 $rec95
